$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the FilesTab query (B4) with the new ICDC query that adds
# File Type, renames Format -> File Format, and adds Breed.
$newFilesQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nWHERE s.clinical_study_designation IN ['UBC01']`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(f.file_type, '') AS ``File Type``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``File Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(demo.breed,'') AS Breed , `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"
$ws.Range("B4").Value = $newFilesQuery

# The new query text takes more lines when wrapped, so row 4 grows to
# match row 3's height.
$ws.Rows.Item(4).RowHeight = 232

# Update the view: scroll down to row 4 and move the active selection to B4.
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("B4").Select()
